$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data-input")

# Rename header tag keys (row 1)
$ws.Range("E1").Value = "artist_name:"
$ws.Range("F1").Value = "album_title:"
$ws.Range("G1").Value = "track_title:"
$ws.Range("H1").Value = "track_number:"

# Rename genre consistency tag values in column K (genre_in_dict:)
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 11)
    $val = $cell.Value2
    if ($val -eq "ARTIST_to_GENRE_OK") {
        $cell.Value = "GENRE_OK"
    } elseif ($val -eq "INCONSISTENT_GENRE") {
        $cell.Value = "INCONSISTENT"
    }
}

# Adjust column widths on data-input sheet
# (runtime quantizes ColumnWidth to 1/6-character steps, so these inputs
# land on the closest achievable stored width to the target 15.7109375 / 16.7109375)
$ws.Columns.Item(8).ColumnWidth = 14.8
$ws.Columns.Item(11).ColumnWidth = 15.8
